$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 21951
$ws.Range("I106").Value = 34918.332
$ws.Range("K106").Value = 34918.332
$ws.Range("M106").Value = -34287.332
$ws.Range("H113").Value = 1991.24
$ws.Range("I113").Value = 1780.6875
$ws.Range("J113").Value = 2365.5557
$ws.Range("K113").Value = 1780.6875
$ws.Range("L113").Value = 2365.5557
$ws.Range("M113").Value = 1473.3125
$ws.Range("N113").Value = -8873.555700000001
$ws.Range("H138").Value = 10754857
$ws.Range("J138").Value = 3584
$ws.Range("L138").Value = 10752
$ws.Range("N138").Value = -21032

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 11906325
$ws.Range("I61").Value = 12821986
$ws.Range("J61").Value = 2743.3333
$ws.Range("K61").Value = 12821986
$ws.Range("L61").Value = 2743.3333
$ws.Range("M61").Value = -12821774
$ws.Range("N61").Value = -3167.3333
$ws.Range("H122").Value = 14144.5
$ws.Range("I122").Value = 14911.363
$ws.Range("J122").Value = 11332.667
$ws.Range("K122").Value = 44734.089
$ws.Range("L122").Value = 33998.001
$ws.Range("M122").Value = -42284.089
$ws.Range("N122").Value = -38898.001
$ws.Range("H136").Value = 11906325
$ws.Range("I136").Value = 12821986
$ws.Range("J136").Value = 2743.3333
$ws.Range("K136").Value = 38465958
$ws.Range("L136").Value = 8229.999899999999
$ws.Range("M136").Value = -38463408
$ws.Range("N136").Value = -13329.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 10500
$ws.Range("I25").Value = 1000
$ws.Range("J25").Value = 20000
$ws.Range("K25").Value = 1000
$ws.Range("L25").Value = 20000
$ws.Range("M25").Value = -826
$ws.Range("N25").Value = -20348
$ws.Range("H31").Value = 7941439
$ws.Range("I31").Value = 7467.0527
$ws.Range("J31").Value = 14495590
$ws.Range("K31").Value = 7467.0527
$ws.Range("L31").Value = 14495590
$ws.Range("M31").Value = -7172.0527
$ws.Range("N31").Value = -14496180
$ws.Range("H34").Value = 7941439
$ws.Range("I34").Value = 7467.0527
$ws.Range("J34").Value = 14495590
$ws.Range("K34").Value = 7467.0527
$ws.Range("L34").Value = 14495590
$ws.Range("M34").Value = -7265.0527
$ws.Range("N34").Value = -14495994
$ws.Range("H58").Value = 1671.6471
$ws.Range("I58").Value = 746.4
$ws.Range("J58").Value = 2993.4285
$ws.Range("K58").Value = 746.4
$ws.Range("L58").Value = 2993.4285
$ws.Range("M58").Value = -543.4
$ws.Range("N58").Value = -3399.4285
$ws.Range("H99").Value = 1961.3889
$ws.Range("I99").Value = 1933.75
$ws.Range("K99").Value = 1933.75
$ws.Range("M99").Value = -435.75
$ws.Range("H109").Value = 24933.334
$ws.Range("J109").Value = 25550
$ws.Range("L109").Value = 25550
$ws.Range("N109").Value = -27630
$ws.Range("H122").Value = 1611.7778
$ws.Range("I122").Value = 1527.4286
$ws.Range("K122").Value = 4582.2858
$ws.Range("M122").Value = -2132.2858
$ws.Range("H126").Value = 1961.3889
$ws.Range("I126").Value = 1933.75
$ws.Range("K126").Value = 5801.25
$ws.Range("M126").Value = -3331.25
$ws.Range("H132").Value = 16668312
$ws.Range("I132").Value = 19232136
$ws.Range("J132").Value = 3457.5
$ws.Range("K132").Value = 57696408
$ws.Range("L132").Value = 10372.5
$ws.Range("M132").Value = -57693878
$ws.Range("N132").Value = -15432.5
$ws.Range("H136").Value = 1671.6471
$ws.Range("I136").Value = 746.4
$ws.Range("J136").Value = 2993.4285
$ws.Range("K136").Value = 2239.2
$ws.Range("L136").Value = 8980.2855
$ws.Range("M136").Value = 310.8000000000002
$ws.Range("N136").Value = -14080.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 1222
$ws.Range("I114").Value = 374.875
$ws.Range("J114").Value = 1578.6842
$ws.Range("K114").Value = 1124.625
$ws.Range("L114").Value = 4736.0526
$ws.Range("M114").Value = 2129.375
$ws.Range("N114").Value = -11244.0526
$ws.Range("H115").Value = 4943.1665
$ws.Range("I115").Value = 1250
$ws.Range("J115").Value = 5681.8
$ws.Range("K115").Value = 3750
$ws.Range("L115").Value = 17045.4
$ws.Range("M115").Value = -2575
$ws.Range("N115").Value = -19395.4
$ws.Range("H121").Value = 617.9
$ws.Range("I121").Value = 344.54544
$ws.Range("J121").Value = 952
$ws.Range("K121").Value = 1033.63632
$ws.Range("L121").Value = 2856
$ws.Range("M121").Value = 276.3636799999999
$ws.Range("N121").Value = -5476

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 14652.111
$ws.Range("I57").Value = 3974.75
$ws.Range("J57").Value = 23194
$ws.Range("K57").Value = 3974.75
$ws.Range("L57").Value = 23194
$ws.Range("M57").Value = -3154.75
$ws.Range("N57").Value = -24834
$ws.Range("H102").Value = 3762.3809
$ws.Range("I102").Value = 5400.846
$ws.Range("J102").Value = 1099.875
$ws.Range("K102").Value = 5400.846
$ws.Range("L102").Value = 1099.875
$ws.Range("M102").Value = -3778.846
$ws.Range("N102").Value = -4343.875
$ws.Range("H122").Value = 4764133.5
$ws.Range("I122").Value = 6668217.5
$ws.Range("K122").Value = 20004652.5
$ws.Range("M122").Value = -20002202.5
$ws.Range("H126").Value = 5300.143
$ws.Range("I126").Value = 3700.5
$ws.Range("J126").Value = 5940
$ws.Range("K126").Value = 11101.5
$ws.Range("L126").Value = 17820
$ws.Range("M126").Value = -8631.5
$ws.Range("N126").Value = -22760

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5806.4116
$ws.Range("I40").Value = 9100.571
$ws.Range("K40").Value = 9100.571
$ws.Range("M40").Value = -8964.571
$ws.Range("H100").Value = 2658.0667
$ws.Range("I100").Value = 2410.111
$ws.Range("K100").Value = 2410.111
$ws.Range("M100").Value = -1869.111
$ws.Range("H122").Value = 5811.846
$ws.Range("I122").Value = 5859
$ws.Range("J122").Value = 5736.4
$ws.Range("K122").Value = 17577
$ws.Range("L122").Value = 17209.2
$ws.Range("M122").Value = -15127
$ws.Range("N122").Value = -22109.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 15000
$ws.Range("J103").Value = 15000
$ws.Range("L103").Value = 15000
$ws.Range("N103").Value = -17344
$ws.Range("H107").Value = 381.83334
$ws.Range("I107").Value = 286.5
$ws.Range("J107").Value = 572.5
$ws.Range("K107").Value = 859.5
$ws.Range("L107").Value = 1717.5
$ws.Range("M107").Value = 1060.5
$ws.Range("N107").Value = -5557.5
$ws.Range("H109").Value = 33377
$ws.Range("J109").Value = 33377
$ws.Range("L109").Value = 33377
$ws.Range("N109").Value = -36151
$ws.Range("H113").Value = 1328.7142
$ws.Range("I113").Value = 220.4
$ws.Range("J113").Value = 1944.4445
$ws.Range("K113").Value = 661.2
$ws.Range("L113").Value = 5833.333500000001
$ws.Range("M113").Value = 1508.8
$ws.Range("N113").Value = -10173.3335
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H115").Value = 34619.75
$ws.Range("J115").Value = 34619.75
$ws.Range("L115").Value = 34619.75
$ws.Range("N115").Value = -37753.75
$ws.Range("H118").Value = 48000
$ws.Range("J118").Value = 48000
$ws.Range("L118").Value = 48000
$ws.Range("N118").Value = -51314
$ws.Range("H136").Value = 1859.3889
$ws.Range("I136").Value = 1798.091
$ws.Range("J136").Value = 1955.7142
$ws.Range("K136").Value = 5394.272999999999
$ws.Range("L136").Value = 5867.142599999999
$ws.Range("M136").Value = -2844.272999999999
$ws.Range("N136").Value = -10967.1426
